# Add a "Total" column to the QuickJobTimeReport sheet (sub-total / grand-total column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G holds the "Total" header; give it a sensible width like the
# other data columns.
$ws.Columns("G").ColumnWidth = 9.8

# Write the new header text first ...
$ws.Range("G5").Value = "Total"

# ... then pick up the shaded "Category" style used across row 5 (same
# font/fill/border treatment as A5/B5/E5/F5) and apply it to the new cell
# so the Total header matches the rest of the row.
$ws.Range("A5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection on the newly added cell, matching where the user
# would land after typing the new header.
$ws.Range("G5").Select()
